$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text (string) data type, matching the
# original inline-string cells, instead of Excel auto-coercing numeric-
# looking text into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "62.574.95"
$ws.Range("D3").Value = "2.436.97"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "568.29"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "145.09"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "2.432.90"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "26.74"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").Value = "2.876.62"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "62.468.04"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.434.93"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "7.27"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").Value = "326.49"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "2.05"
$ws.Range("E23").Value = "  +12.03%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "65.01"
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("D26").Value = "614.17"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "8.75"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").Value = "2.558.94"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").Value = "5.10"
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "18.70"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "5.31"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").Value = "145.01"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").Value = "147.55"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "3.74"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "20.81"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "0.0529"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").Value = "0.595"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("E51").Value = "  -0.67%  "
